$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1875.1923
$ws.Range("I15").Value = 1875.1923
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 5625.5769
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -5456.5769

$ws.Range("H88").Value = 2708.25
$ws.Range("I88").Value = 2425
$ws.Range("J88").Value = 2802.6667
$ws.Range("K88").Value = 2425
$ws.Range("L88").Value = 2802.6667
$ws.Range("M88").Value = -2019
$ws.Range("N88").Value = -3614.6667

$ws.Range("H91").Value = 2708.25
$ws.Range("I91").Value = 2425
$ws.Range("J91").Value = 2802.6667
$ws.Range("K91").Value = 2425
$ws.Range("L91").Value = 2802.6667
$ws.Range("M91").Value = -1021
$ws.Range("N91").Value = -5610.6667

$ws.Range("H92").Value = 250001800
$ws.Range("I92").Value = 250001800
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 250001800
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -250000552

$ws.Range("H129").Value = 919.9149
$ws.Range("I129").Value = 384.92307
$ws.Range("J129").Value = 1124.4706
$ws.Range("K129").Value = 1154.76921
$ws.Range("L129").Value = 3373.4118
$ws.Range("M129").Value = 3845.23079
$ws.Range("N129").Value = -13373.4118

$ws.Range("H131").Value = 1759
$ws.Range("I131").Value = 265
$ws.Range("J131").Value = 4000
$ws.Range("K131").Value = 795
$ws.Range("L131").Value = 12000
$ws.Range("M131").Value = 4245
$ws.Range("N131").Value = -22080

$ws.Range("H137").Value = 1493.375
$ws.Range("I137").Value = 1353.8077
$ws.Range("J137").Value = 2098.1667
$ws.Range("K137").Value = 4061.4231
$ws.Range("L137").Value = 6294.500100000001
$ws.Range("M137").Value = -1511.4231

$ws.Range("H138").Value = 2196.987
$ws.Range("I138").Value = 1928.0555
$ws.Range("J138").Value = 2277.6667
$ws.Range("K138").Value = 5784.166499999999
$ws.Range("L138").Value = 6833.000100000001
$ws.Range("M138").Value = -644.1664999999994
$ws.Range("N138").Value = -17113.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6554.3477
$ws.Range("I32").Value = 4704.5454
$ws.Range("J32").Value = 47250
$ws.Range("K32").Value = 4704.5454
$ws.Range("L32").Value = 47250
$ws.Range("M32").Value = -4417.5454

$ws.Range("H45").Value = 2028.2122
$ws.Range("I45").Value = 2669.5386
$ws.Range("J45").Value = 1611.35
$ws.Range("K45").Value = 2669.5386
$ws.Range("L45").Value = 1611.35
$ws.Range("M45").Value = -2292.5386
$ws.Range("N45").Value = -2365.35

$ws.Range("H74").Value = 52632532
$ws.Range("I74").Value = 71429256
$ws.Range("J74").Value = 1700
$ws.Range("K74").Value = 71429256
$ws.Range("L74").Value = 1700
$ws.Range("M74").Value = -71428382

$ws.Range("H77").Value = 52632532
$ws.Range("I77").Value = 71429256
$ws.Range("J77").Value = 1700
$ws.Range("K77").Value = 357146280
$ws.Range("L77").Value = 8500
$ws.Range("M77").Value = -357141912

$ws.Range("H122").Value = 2957.6316
$ws.Range("I122").Value = 2528.4285
$ws.Range("J122").Value = 4159.4
$ws.Range("K122").Value = 7585.2855
$ws.Range("L122").Value = 12478.2
$ws.Range("M122").Value = -5135.2855

$ws.Range("H132").Value = 11830.708
$ws.Range("I132").Value = 1239.9048
$ws.Range("J132").Value = 85966.336
$ws.Range("K132").Value = 3719.7144
$ws.Range("L132").Value = 257899.008
$ws.Range("M132").Value = -1189.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1281.2222
$ws.Range("I16").Value = 1281.2222
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1281.2222
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -994.2221999999999
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 15792.654
$ws.Range("I31").Value = 30218.273
$ws.Range("J31").Value = 5213.8667
$ws.Range("K31").Value = 30218.273
$ws.Range("L31").Value = 5213.8667
$ws.Range("M31").Value = -29923.273

$ws.Range("H34").Value = 15792.654
$ws.Range("I34").Value = 30218.273
$ws.Range("J34").Value = 5213.8667
$ws.Range("K34").Value = 30218.273
$ws.Range("L34").Value = 5213.8667
$ws.Range("M34").Value = -30016.273

$ws.Range("H113").Value = 1281.2222
$ws.Range("I113").Value = 1281.2222
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1281.2222
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 888.7778000000001
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 13389.889
$ws.Range("I132").Value = 19134.414
$ws.Range("J132").Value = 2977.9375
$ws.Range("K132").Value = 57403.242
$ws.Range("L132").Value = 8933.8125
$ws.Range("M132").Value = -54873.242
$ws.Range("N132").Value = -13993.8125

$ws.Range("H134").Value = 1339.7273
$ws.Range("I134").Value = 1215.375
$ws.Range("J134").Value = 1671.3334
$ws.Range("K134").Value = 3646.125
$ws.Range("L134").Value = 5014.0002
$ws.Range("M134").Value = -1111.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 784.48
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 784.48
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2353.44
$ws.Range("N131").Value = -12433.44
$ws.Range("M131").ClearContents()

$ws.Range("H132").Value = 1352.4166
$ws.Range("I132").Value = 475
$ws.Range("J132").Value = 1527.9
$ws.Range("K132").Value = 4275
$ws.Range("L132").Value = 13751.1
$ws.Range("M132").Value = -1745
$ws.Range("N132").Value = -18811.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H122").Value = 111112020
$ws.Range("I122").Value = 33334120
$ws.Range("J122").Value = 500001500
$ws.Range("K122").Value = 100002360
$ws.Range("L122").Value = 1500004500
$ws.Range("M122").Value = -99999910

$ws.Range("H130").Value = 39265.316
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 39265.316
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 39265.316
$ws.Range("N130").Value = -49305.316

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4906.7
$ws.Range("I22").Value = 4057.2856
$ws.Range("J22").Value = 6888.6665
$ws.Range("K22").Value = 4057.2856
$ws.Range("L22").Value = 6888.6665
$ws.Range("M22").Value = -3762.2856
$ws.Range("N22").Value = -7478.6665

$ws.Range("H27").Value = 4906.7
$ws.Range("I27").Value = 4057.2856
$ws.Range("J27").Value = 6888.6665
$ws.Range("K27").Value = 4057.2856
$ws.Range("L27").Value = 6888.6665
$ws.Range("M27").Value = -3950.2856
$ws.Range("N27").Value = -7102.6665

$ws.Range("H46").Value = 984.0769
$ws.Range("I46").Value = 999.2222
$ws.Range("J46").Value = 950
$ws.Range("K46").Value = 999.2222
$ws.Range("L46").Value = 950
$ws.Range("M46").Value = -811.2222
$ws.Range("N46").Value = -1326

$ws.Range("H55").Value = 85.30768999999999
$ws.Range("I55").Value = 87
$ws.Range("J55").Value = 84.25
$ws.Range("K55").Value = 87
$ws.Range("L55").Value = 84.25
$ws.Range("M55").Value = 86
$ws.Range("N55").Value = -430.25

$ws.Range("H93").Value = 3741.8572
$ws.Range("I93").Value = 3364.8333
$ws.Range("J93").Value = 6004
$ws.Range("K93").Value = 3364.8333
$ws.Range("L93").Value = 6004
$ws.Range("M93").Value = -2116.8333
$ws.Range("N93").Value = -8500

$ws.Range("H122").Value = 1034810.4
$ws.Range("I122").Value = 2181343.8
$ws.Range("J122").Value = 2930.4
$ws.Range("K122").Value = 6544031.399999999
$ws.Range("L122").Value = 8791.200000000001
$ws.Range("M122").Value = -6541581.399999999
$ws.Range("N122").Value = -13691.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 23728.666
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 23728.666
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 23728.666
$ws.Range("N46").Value = -24190.666

$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H132").Value = 779.8214
$ws.Range("I132").Value = 538.15216
$ws.Range("J132").Value = 1891.5
$ws.Range("K132").Value = 1614.45648
$ws.Range("L132").Value = 5674.5
$ws.Range("M132").Value = 915.5435200000002

$ws.Range("H134").Value = 23728.666
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 23728.666
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 71185.99800000001
$ws.Range("N134").Value = -76255.99800000001
